$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "CH-S09FTXF-NG"
$ws.Cells.Item(16, 3).Value = 610
$ws.Cells.Item(16, 4).Value = 762.5

# Row 17
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = 3
$ws.Cells.Item(17, 3).Value = "CH-S12FTXF-NG"
$ws.Cells.Item(17, 4).Value = 661
$ws.Cells.Item(17, 5).Value = 826.25

# Row 18
$ws.Cells.Item(18, 1).Value = "CH-S09FTXF-NG"
$ws.Cells.Item(18, 2).Value = 610
$ws.Cells.Item(18, 3).Value = 762.5

# Row 19
$ws.Cells.Item(19, 1).Value = "CH-S12FTXF-NG"
$ws.Cells.Item(19, 2).Value = 661
$ws.Cells.Item(19, 3).Value = 826.25

# Row 20
$ws.Cells.Item(20, 1).Value = "CH-S18FTXF-NG"
$ws.Cells.Item(20, 2).Value = 950
$ws.Cells.Item(20, 3).Value = 1187.5

# Row 21
$ws.Cells.Item(21, 1).Value = "CH-S09FTXF-NG"
$ws.Cells.Item(21, 2).Value = 1408
$ws.Cells.Item(21, 3).Value = 1760

# Row 22
$ws.Cells.Item(22, 1).Value = "CH-S18FTXLA2-NG"
$ws.Cells.Item(22, 2).Value = 724
$ws.Cells.Item(22, 3).Value = 905

# Row 23
$ws.Cells.Item(23, 1).Value = "CH-S12FTXLA2-NG"
$ws.Cells.Item(23, 2).Value = 790
$ws.Cells.Item(23, 3).Value = 987.5

# Row 24
$ws.Cells.Item(24, 1).Value = "CH-S24FTXLA2-NG"
$ws.Cells.Item(24, 2).Value = 1269
$ws.Cells.Item(24, 3).Value = 1586.25
